$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-5 (columns A:AH / 1:34)
$rowData = @(
    @(45087.50694444445, 17.776, 11.78, 4.002, 37.831, 29.945, 13.989, 43.56, 21.524, 8.901999999999999, 13.298, 14.882, 15.457, 4.465, 13.911, 19.348, 12.001, 3.404, 2.232, 204.374, 38.641, 12.84, 25.287, 12.826, 3.128, 22.227, 11.342, 10.33, 12.137, 15.411, 3.318, 38.832, 6.978, 16.052),
    @(45087.51388888889, 7.687, 5.053, 1.519, 16.544, 12.856, 6.05, 25.796, 9.308, 3.86, 5.522, 6.592, 6.799, 1.936, 6.015, 8.363, 5.424, 1.463, 0.822, 84.25700000000001, 16.981, 5.552, 10.986, 5.59, 1.437, 12.301, 4.904, 4.584, 5.358, 6.788, 1.247, 23.835, 2.939, 6.943),
    @(45087.52083333334, 22.58, 16.496, 1.528, 49.043, 39.94, 17.77, 65.62, 27.341, 12.057, 17.817, 19.671, 20.697, 5.675, 17.67, 25.073, 14.977, 1.093, 1.009, 261.637, 49.345, 16.31, 33.059, 17.308, 2.75, 32.379, 14.407, 12.836, 15.081, 20.596, 0.766, 59.537, 9.132999999999999, 20.392),
    @(45087.52777777778, 4.32, 2.9, 0.7, 9.359999999999999, 7.25, 3.4, 19.66, 5.24, 2.26, 3.1, 3.76, 3.89, 1.1, 3.38, 4.8, 3.08, 0.71, 0.36, 44.25, 9.83, 3.12, 6.44, 3.24, 0.77, 8.92, 2.76, 2.58, 3.01, 3.89, 0.55, 18.44, 1.65, 3.91)
)

for ($i = 0; $i -lt $rowData.Length; $i++) {
    $r = $i + 2
    $vals = $rowData[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $vals[$j]
    }
}

# Remove the old row 6 (data now only spans rows 1-5)
$ws.Rows.Item(6).Delete()

# Widen a subset of columns from 7 to 8 characters
# (ColumnWidth uses a different unit than the stored XML "width"; offset of 0.83 observed)
$wideCols = @(2, 3, 7, 10, 11, 13, 17, 24, 27, 28, 29)
foreach ($c in $wideCols) {
    $ws.Columns.Item($c).ColumnWidth = 7.17
}
